# Adds a new paragraph "Which is your best part of this class? Nikita"
# right after the existing "How many assignments do you have due this
# week? Nikita" paragraph, at the very end of the document.
#
# The trailing "_GoBack" bookmark (wrapping the previously-last
# paragraph) must end up wrapping the new, now-last paragraph instead
# -- exactly like Word does when you place the caret at the end of the
# document and type a new line. We achieve that by inserting the new
# text and the paragraph break directly at the bookmark's collapsed
# location: Word (and this host) keeps a zero-width bookmark anchored
# immediately before text that is inserted right at its position, so it
# rides along to the end of the newly typed content.

$d = $word.ActiveDocument

$anchorText = "How many assignments do you have due this week? Nikita"
$newText = "Which is your best part of this class? Nikita"

# Locate the end of the anchor paragraph's text irrespective of exactly
# where the "_GoBack" bookmark happens to live.
$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph text."
}
$insertPos = $searchRange.End

# 1) Insert the new paragraph's text right after the anchor text (still
#    inside the same paragraph for now).
$textRange = $d.Range($insertPos, $insertPos)
$textRange.InsertAfter($newText)

# 2) Split that paragraph in two by inserting a paragraph mark at the
#    original boundary, so the new text becomes its own paragraph.
$breakRange = $d.Range($insertPos, $insertPos)
$breakRange.InsertAfter([char]13)
